$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Copy formatting (style indices) from row 26 (A:C) onto row 27 (A:C) so the
# newly-populated cells pick up the same styles used throughout the table.
$ws.Range("A26:C26").Copy()
$ws.Range("A27:C27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 27 - PUBLONS025
$ws.Range("A27").Value = 'PUBLONS025'
$ws.Range("B27").Value = 'OPQA-5893||OPQA-5894||OPQA-5899'
$ws.Range("C27").Value = 'Verify the linking model should display when user sign in using facebook and enters matching account from the login page||Verify the User is redirected to Publon''s Homepage when User enters correct credentials for matching facebook account and Also verify user can see two connected  account on account setting page one is facebook and another is Steam account.||Verify User should not challenged linking model when Facebook is already linked with STeAM.'
$ws.Range("D27").Value = 'Y'
$ws.Rows.Item(27).RowHeight = 60

# Row 28 - PUBLONS026
$ws.Range("A28").Value = 'PUBLONS026'
$ws.Range("B28").Value = 'OPQA-5897||OPQA-5898||OPQA-5999'
$ws.Range("C28").Value = 'Verify the linking model should display when user sign in using Gmail and enters matching account from the login page||Verify the User is redirected to Publon''s Homepage when User enters correct credentials for matching Gmail account and Also verify user can see two connected  account on account setting page one is Gmail and another is Steam account.||Verify User should not challenged linking model when Gmail is already linked with STeAM.'
$ws.Range("D28").Value = 'Y'
$ws.Rows.Item(28).RowHeight = 60

# Update the sheet view to reflect the scrolled/selected state after entry
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C33").Select()

